$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-5 (D=date serial, J=Volumen, K/L/M=Precio min/max/prom, P=Precio $/Kg)
$targets = @{
    2 = @{ D = 44348; J = 20; K = 10000; L = 10000; M = 10000; P = 1000 }
    3 = @{ D = 44463; J = 25; K = 12000; L = 12000; M = 12000; P = 1200 }
    4 = @{ D = 44473; J = 25; K = 11000; L = 11000; M = 11000; P = 1100 }
    5 = @{ D = 44469; J = 20; K = 12000; L = 12000; M = 12000; P = 1200 }
}

foreach ($row in $targets.Keys) {
    $vals = $targets[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
